# Applies the "New simulation files for schemes report" edit:
#  - HKL label rows (B column) get new/renamed scheme names, and the
#    [h,k,l] / pairing headers in row 2 are re-ordered.
#  - 10 new schemes (rows 20-29) are appended.
#  - The now-unused trailing columns U:AD (duplicate 19-28 index row +
#    duplicate header row values) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the now-unused columns U:AD in rows 1 and 2 (table shrinks to
#    20 columns / A:T) so the sheet dimension becomes A1:T29.
# ---------------------------------------------------------------------
$ws.Range("U1:AD2").Clear()

# ---------------------------------------------------------------------
# 2. Prepare rows 20-29 (new schemes) with the same look as row 19:
#    bold/bordered index in column A, plain value cells elsewhere.
# ---------------------------------------------------------------------
$ws.Range("A19:T19").Copy()
$ws.Range("A20:T29").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Row 1 (column index strip 0..18) is unchanged other than losing the
#    trailing 19..28 values handled by the Clear() above.
# ---------------------------------------------------------------------
$row1 = New-Object 'object[,]' 1,19
for ($i = 0; $i -le 18; $i++) { $row1[0,$i] = $i }
$ws.Range("B1:T1").Value2 = $row1

# ---------------------------------------------------------------------
# 4. Row 2: A2 stays 0/HKL; C2:T2 become the [h,k,l] reflections followed
#    by the pairing-count labels (re-ordered vs. the old layout).
# ---------------------------------------------------------------------
$ws.Range("A2").Value2 = 0
$ws.Range("B2").Value2 = "HKL"
$row2 = New-Object 'object[,]' 1,18
$row2Vals = @(
    "[4, 0, 0]", "[2, 1, 1]", "[2, 2, 0]", "[2, 0, 0]", "[2, 2, 2]", "[3, 1, 0]",
    "[1, 1, 0]", "[3, 2, 1]",
    "1Pair-A", "1Pair-B", "2Pairs-A", "2Pairs-B", "3Pairs-A", "3Pairs-B",
    "3Pairs-C", "4Pairs", "5A4F", "MaxUnique"
)
for ($i = 0; $i -lt $row2Vals.Length; $i++) { $row2[0,$i] = $row2Vals[$i] }
$ws.Range("C2:T2").Value2 = $row2

# ---------------------------------------------------------------------
# 5. Column B (scheme names) for rows 3-29 -- the C:T "1" flag cells in
#    rows 3-19 are untouched, rows 20-29 already got 1's pasted in step 2
#    and are overwritten below with their real index/name in column A/B.
# ---------------------------------------------------------------------
$schemeNames = @(
    "Spiral5",                    # row 3  (index 1)
    "RotRing OmegaMax-90",        # row 4  (index 2)
    "Equal Angle",                # row 5  (index 3)
    "Tilt Rotate",                # row 6  (index 4)
    "CLR",                        # row 7  (index 5)
    "Rizzie Hex",                 # row 8  (index 6)
    "Thomas Hex",                 # row 9  (index 7)
    "Tilt Rotate_Partial",        # row 10 (index 8)
    "RotRing OmegaMax-60",        # row 11 (index 9)
    "Equal Angle_Partial",        # row 12 (index 10)
    "Rizzie Hex_Partial",         # row 13 (index 11)
    "ND Single",                  # row 14 (index 12)
    "RD Single",                  # row 15 (index 13)
    "TD Single",                  # row 16 (index 14)
    "Morris Single",              # row 17 (index 15)
    "Ring Perpendicular to ND",   # row 18 (index 16)
    "Ring Perpendicular to RD",   # row 19 (index 17)
    "Ring Perpendicular to TD",   # row 20 (index 18)
    "OffsetFTD",                  # row 21 (index 19)
    "OffsetATD",                  # row 22 (index 20)
    "OffsetF45",                  # row 23 (index 21)
    "OffsetA45",                  # row 24 (index 22)
    "OffsetFRD",                  # row 25 (index 23)
    "OffsetARD",                  # row 26 (index 24)
    "Gaussian Quadrature",        # row 27 (index 25)
    "Michael-CCHex",              # row 28 (index 26)
    "Michael-SNHex"               # row 29 (index 27)
)

for ($i = 0; $i -lt $schemeNames.Length; $i++) {
    $r = 3 + $i
    $ws.Range("A$r").Value2 = $i + 1
    $ws.Range("B$r").Value2 = $schemeNames[$i]
}

# ---------------------------------------------------------------------
# 6. Fill the "1" indicator cells (C:T) for the 10 brand-new rows 20-29.
# ---------------------------------------------------------------------
$onesRow = New-Object 'object[,]' 1,18
for ($i = 0; $i -lt 18; $i++) { $onesRow[0,$i] = 1 }
for ($r = 20; $r -le 29; $r++) {
    $ws.Range("C${r}:T${r}").Value2 = $onesRow
}
